$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2041586444627012
$ws.Range("C2").Value = -0.317315650508697
$ws.Range("D2").Value = 1.506101047285106

$ws.Range("B3").Value = -0.20815504766357
$ws.Range("C3").Value = 1.460937717351997

$ws.Range("B4").Value = 1.167899071513181

$ws.Range("B5").Value = 1.713367536971843
$ws.Range("C5").Value = -2.474187431820313
$ws.Range("D5").Value = 0.4111563609322812
$ws.Range("E5").Value = 0.3036778841160683

$ws.Range("B6").Value = -1.934752120961668
$ws.Range("C6").Value = 0.3549748406297492
$ws.Range("D6").Value = 0.2564759890460406

$ws.Range("B7").Value = -0.9091272736756697
$ws.Range("C7").Value = 0.8379364208866271

$ws.Range("B8").Value = 0.3124575324117352

$ws.Range("B9").Value = 0.5586933435139947
$ws.Range("C9").Value = 0.1779189921522184
$ws.Range("D9").Value = 0.3037985701168242
$ws.Range("E9").Value = 0.295512446150037

$ws.Range("B10").Value = 0.4096934484336369
$ws.Range("C10").Value = 0.28795472432474
$ws.Range("D10").Value = 0.3209312923645784

$ws.Range("B11").Value = 0.4872844834843318
$ws.Range("C11").Value = 0.2873668796881647

$ws.Range("B12").Value = 0.4220310903795553

$ws.Range("B13").Value = 0.5981638006715656
$ws.Range("C13").Value = -0.05550181425120998
$ws.Range("D13").Value = 0.04865461862192938

$ws.Range("B14").Value = 0.168998243652269
$ws.Range("C14").Value = 0.08624861542179052

$ws.Range("B15").Value = 0.1585083094968939
